# Workbook edit: "Added all tables from the spreadsheet and modified some
# date names (date is a keyword)"
#
# This database-diagram worksheet lays out tables in two side-by-side
# columns (A/B and D/E). Before this edit, the EMPLOYEE table sat in A27:B31
# and a separate SCHEDULE table sat below it in A35:B44. This change moves
# the SCHEDULE table so it sits next to EMPLOYEE (in D27:E36, mirroring the
# A/B table-header/field layout), renames its "date" field to
# "schedule_date", and renames the CARE table's "date" field (D15/E15) to
# "care_date" (fixing its type to "date NOT NULL" at the same time).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the old SCHEDULE table's tail rows (A40:B44 -> care_blue..
#    care_yellow) entirely; they get re-created under D/E further down.
# ---------------------------------------------------------------------
$ws.Range("A40:A44").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2. Rows 37 ("date"/"date NOT NULL") and 38 ("made_by"/...) of the old
#    SCHEDULE table are no longer needed either (their content moves to
#    D29/D30). Deleting them shifts row 39 (the blank bottom-border cells
#    D39:E39) up to row 37, so immediately re-insert two blank rows to
#    push that content back down to row 39 where it belongs.
# ---------------------------------------------------------------------
$ws.Range("A37:A38").EntireRow.Delete()
$ws.Range("A37:A38").EntireRow.Insert()

# ---------------------------------------------------------------------
# 3. The old SCHEDULE header (A35:B35) was a merged cell; unmerge it,
#    then clear the leftover A/B cells of the old table (rows 35-36 still
#    hold "SCHEDULE"/"schedule_id", row 39 still holds "doctor_id").
# ---------------------------------------------------------------------
$ws.Range("A35:B35").UnMerge()
$ws.Range("A35:B36").Clear()
$ws.Range("A39:B39").Clear()

# ---------------------------------------------------------------------
# 4. Write the SCHEDULE table into its new home beside EMPLOYEE
#    (D27:E36). D27/E27 already carry the bold/centered header style
#    from the template, so only the value needs to be set there.
# ---------------------------------------------------------------------
$ws.Range("D27").Value = "SCHEDULE"

$ws.Range("D28").Value = "schedule_id"
$ws.Range("E28").Value = "int PK"

$ws.Range("D29").Value = "schedule_date"
$ws.Range("E29").Value = "date NOT NULL"

$ws.Range("D30").Value = "made_by"
$ws.Range("E30").Value = "int FK (EMPLOYEE.emp_id) NOT NULL"

$ws.Range("D31").Value = "doctor_id"
$ws.Range("E31").Value = "int FK (EMPLOYEE.emp_id) NOT NULL"

$ws.Range("D32").Value = "supervisor_id"
$ws.Range("E32").Value = "int FK (EMPLOYEE.emp_id) NOT NULL"

$ws.Range("D33").Value = "care_red"
$ws.Range("E33").Value = "int FK (EMPLOYEE.emp_id) NOT NULL"

$ws.Range("D34").Value = "care_blue"
$ws.Range("E34").Value = "int FK (EMPLOYEE.emp_id) NOT NULL"

$ws.Range("D35").Value = "care_green"
$ws.Range("E35").Value = "int FK (EMPLOYEE.emp_id) NOT NULL"

$ws.Range("D36").Value = "care_yellow"
$ws.Range("E36").Value = "int FK (EMPLOYEE.emp_id) NOT NULL"

# ---------------------------------------------------------------------
# 5. CARE table: rename the "date" field to "care_date" ("date" is a SQL
#    keyword) and fix its type to "date NOT NULL" to match the other
#    date columns.
# ---------------------------------------------------------------------
$ws.Range("D15").Value = "care_date"
$ws.Range("E15").Value = "date NOT NULL"

# ---------------------------------------------------------------------
# 6. Restore the author's cursor position.
# ---------------------------------------------------------------------
$ws.Range("E26").Select()
